# Update CDA Logical model metadata (ST.r2b) -- StructureDefinition-ANY.xlsx
#
# Changes applied to the "Metadata" worksheet (sheet1):
#   - Version value bumped: 2.0.0-sd-202406-matchbox-patch -> 2.0.1-sd-202510-matchbox-patch
#   - Date value bumped:   2024-06-19T17:47:42+02:00 -> 2025-10-29T22:15:57+01:00
#   - New "Jurisdiction" property/value row inserted right after "Contact"
#     (before "Description"), pushing the remaining rows down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# --- Version ---------------------------------------------------------
$ws.Range("B3").Value = "2.0.1-sd-202510-matchbox-patch"

# --- Date -------------------------------------------------------------
$ws.Range("B8").Value = "2025-10-29T22:15:57+01:00"

# --- Insert the new "Jurisdiction" row after "Contact" (row 10) ------
# "Contact" is row 10, "Description" is row 11 before the insert, so the
# new row goes in at row 11 and everything from the old row 11 onward
# shifts down by one.
$ws.Rows.Item(11).Insert()

# Match the formatting used by the rest of the data rows (style "2":
# top-aligned, wrapped, bordered) by copying it down from the row that
# used to be row 11 ("Description", now row 12) onto the freshly
# inserted row.
$ws.Range("A12:B12").Copy()
$ws.Range("A11:B11").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("A11").Value = "Jurisdiction"
# Value (B11) is intentionally left blank -- matches the other
# not-yet-populated metadata rows (e.g. Purpose / Copyright).
